$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "Khrystyne- Approaches"
$ws.Range("E6").Value = "Derek- Indifferent"
$ws.Range("F6").Value = "Derek- Runs away"
$ws.Range("D5").Value = "Khrystyne- Tail Twitches"
$ws.Range("E5").Value = "Khrystyne- Tail Flags"

$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Columns("D").ColumnWidth = 22.85546875
$ws.Columns("E").ColumnWidth = 19.28515625
$ws.Columns("F").ColumnWidth = 16.85546875

$ws.Range("D1").Select()
